$wb = $excel.ActiveWorkbook

$assay = $wb.Worksheets.Item("assay")
$files = $wb.Worksheets.Item("files")

$tiffNames = @(
  "fov_0_H_0_C_0.tiff",
  "fov_0_H_0_C_1.tiff",
  "fov_0_H_0_C_2.tiff",
  "fov_0_H_0_C_3.tiff",
  "fov_0_H_1_C_0.tiff",
  "fov_0_H_1_C_1.tiff",
  "fov_0_H_1_C_2.tiff",
  "fov_0_H_1_C_3.tiff",
  "fov_0_H_2_C_0.tiff",
  "fov_0_H_2_C_1.tiff",
  "fov_0_H_2_C_2.tiff",
  "fov_0_H_2_C_3.tiff",
  "fov_0_H_3_C_0.tiff",
  "fov_0_H_3_C_1.tiff",
  "fov_0_H_3_C_2.tiff",
  "fov_0_H_3_C_3.tiff",
  "fov_0_dapi.tiff",
  "fov_0_dots.tiff"
)

# First fill in the filenames for the 18 image files on both the
# assay and files tabs.
for ($i = 0; $i -lt $tiffNames.Length; $i++) {
    $row = 3 + $i
    $name = $tiffNames[$i]
    $assay.Cells.Item($row, 10).Value = $name
    $files.Cells.Item($row, 1).Value = $name
}

# Then fill in the format column for those same rows.
for ($i = 0; $i -lt $tiffNames.Length; $i++) {
    $row = 3 + $i
    $files.Cells.Item($row, 2).Value = ".tiff"
}

# Finally add the org.json metadata file as the last row.
$lastRow = 3 + $tiffNames.Length
$assay.Cells.Item($lastRow, 10).Value = "org.json"
$files.Cells.Item($lastRow, 1).Value = "org.json"
$files.Cells.Item($lastRow, 2).Value = ".json"

# Selection / active sheet changes
$assay.Range("I6").Select()
$files.Select()
$files.Range("A26").Select()
